$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.807.00'
$ws.Range("E2").Value = '  -1.25%  '

# Row 3
$ws.Range("D3").Value = '1.854.52'
$ws.Range("E3").Value = '  -0.73%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '''303.83'
$ws.Range("E5").Value = '  -1.04%  '

# Row 6
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.05%  '

# Row 7
$ws.Range("D7").Value = '''0.5043'
$ws.Range("E7").Value = '  -2.16%  '

# Row 8
$ws.Range("D8").Value = '''0.3648'
$ws.Range("E8").Value = '  -2.81%  '

# Row 9
$ws.Range("D9").Value = '''0.07153'
$ws.Range("E9").Value = '  -0.28%  '

# Row 10
$ws.Range("D10").Value = '''0.8896'
$ws.Range("E10").Value = '  +0.28%  '

# Row 11
$ws.Range("D11").Value = '''20.61'
$ws.Range("E11").Value = '  -0.40%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.860.15'
$ws.Range("E12").Value = '  -0.26%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '''0.07477'
$ws.Range("E13").Value = '  -1.27%  '

# Row 14
$ws.Range("D14").Value = '''91.95'
$ws.Range("E14").Value = '  +2.84%  '

# Row 15
$ws.Range("D15").Value = '''5.220'
$ws.Range("E15").Value = '  -2.16%  '

# Row 16
$ws.Range("D16").Value = '''1.002'
$ws.Range("E16").Value = '  +0.01%  '

# Row 17
$ws.Range("D17").Value = '''0.000008490'
$ws.Range("E17").Value = '  -0.89%  '

# Row 18
$ws.Range("D18").Value = '''14.02'
$ws.Range("E18").Value = '  -1.14%  '

# Row 19
$ws.Range("E19").Value = '  -0.02%  '

# Row 20
$ws.Range("D20").Value = '26.847.67'
$ws.Range("E20").Value = '  -1.25%  '

# Row 21
$ws.Range("D21").Value = '''5.014'
$ws.Range("E21").Value = '  -0.81%  '

# Row 22
$ws.Range("D22").Value = '2.091.99'
$ws.Range("E22").Value = '  +0.08%  '

# Row 23
$ws.Range("D23").Value = '''10.31'
$ws.Range("E23").Value = '  -2.92%  '

# Row 24
$ws.Range("D24").Value = '''6.431'
$ws.Range("E24").Value = '  -0.71%  '

# Row 25
$ws.Range("D25").Value = '''146.78'
$ws.Range("E25").Value = '  -2.74%  '

# Row 26
$ws.Range("D26").Value = '''1.789'
$ws.Range("E26").Value = '  -3.10%  '

# Row 27
$ws.Range("D27").Value = '''17.78'
$ws.Range("E27").Value = '  -1.30%  '

# Row 28
$ws.Range("D28").Value = '''2.058'
$ws.Range("E28").Value = '  -3.08%  '

# Row 29
$ws.Range("D29").Value = '''112.66'
$ws.Range("E29").Value = '  -0.03%  '

# Row 30
$ws.Range("D30").Value = '''4.622'
$ws.Range("E30").Value = '  -2.91%  '

# Row 31
$ws.Range("D31").Value = '''4.647'
$ws.Range("E31").Value = '  -0.94%  '

# Row 32
$ws.Range("D32").Value = '''0.09196'
$ws.Range("E32").Value = '  +2.18%  '

# Row 33
$ws.Range("D33").Value = '''0.05073'
$ws.Range("E33").Value = '  -1.51%  '

# Row 34
$ws.Range("D34").Value = '''2.992'
$ws.Range("E34").Value = '  -3.47%  '

# Row 35
$ws.Range("D35").Value = '''0.7433'
$ws.Range("E35").Value = '  -1.24%  '

# Row 36
$ws.Range("D36").Value = '''1.142'
$ws.Range("E36").Value = '  -2.48%  '

# Row 37
$ws.Range("D37").Value = '''3.236'
$ws.Range("E37").Value = '  +6.79%  '

# Row 38
$ws.Range("D38").Value = '''2.502'
$ws.Range("E38").Value = '  -1.01%  '

# Row 39
$ws.Range("D39").Value = '''0.01980'
$ws.Range("E39").Value = '  -2.83%  '

# Row 40
$ws.Range("D40").Value = '''1.077'
$ws.Range("E40").Value = '  -0.30%  '

# Row 41
$ws.Range("D41").Value = '''0.5308'
$ws.Range("E41").Value = '  -0.81%  '

# Row 42
$ws.Range("D42").Value = '''119.41'
$ws.Range("E42").Value = '  +4.14%  '

# Row 43
$ws.Range("D43").Value = '''6.459'
$ws.Range("E43").Value = '  -2.73%  '

# Row 44
$ws.Range("D44").Value = '''8.347'
$ws.Range("E44").Value = '  -1.59%  '

# Row 45
$ws.Range("D45").Value = '''0.1452'
$ws.Range("E45").Value = '  -2.10%  '

# Row 46
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '''1.001'
$ws.Range("E46").Value = '  +0.04%  '

# Row 47
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.4628'
$ws.Range("E47").Value = '  -0.85%  '

# Row 48
$ws.Range("D48").Value = '''9.930'
$ws.Range("E48").Value = '  -2.03%  '

# Row 49
$ws.Range("D49").Value = '''1.553'
$ws.Range("E49").Value = '  -1.24%  '

# Row 50
$ws.Range("D50").Value = '''36.88'
$ws.Range("E50").Value = '  +1.11%  '

# Row 51
$ws.Range("D51").Value = '''62.75'
$ws.Range("E51").Value = '  -3.40%  '
